# Rows 9/10, 15-18, 30/31, 32/33 each form a group of records that were
# re-ordered (their field values rotated among the rows of the group) in the
# source export. Re-apply that rotation by writing each row's new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 131066783
$ws.Range("B9").Value = 83089
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 1312
$ws.Range("F9").Value = "Gammelgransskål"
$ws.Range("G9").Value = "Pseudographis pinicola"
$ws.Range("H9").Value = "(Nyl.) Rehm"
$ws.Range("K9").Value = ""
$ws.Range("L9").Value = ""
$ws.Range("M9").Value = ""
$ws.Range("N9").Value = ""
$ws.Range("Q9").Value = 425170
$ws.Range("R9").Value = 6712292
$ws.Range("AC9").Value = ""

# Row 10
$ws.Range("A10").Value = 131066770
$ws.Range("B10").Value = 57884
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("K10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = "äldre spår"
$ws.Range("N10").Value = ""
$ws.Range("Q10").Value = 425323
$ws.Range("R10").Value = 6712206
$ws.Range("AC10").Value = "Ringhack på gran"

# Row 15
$ws.Range("A15").Value = 131066761
$ws.Range("B15").Value = 91771
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 5447
$ws.Range("F15").Value = "Vedticka"
$ws.Range("G15").Value = "Fuscoporia viticola"
$ws.Range("H15").Value = "(Schwein.) Murrill"
$ws.Range("K15").Value = ""
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = ""
$ws.Range("Q15").Value = 425072
$ws.Range("R15").Value = 6712273
$ws.Range("AC15").Value = ""

# Row 16
$ws.Range("A16").Value = 131066782
$ws.Range("B16").Value = 91822
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 1204
$ws.Range("F16").Value = "Gränsticka"
$ws.Range("G16").Value = "Phellopilus nigrolimitatus"
$ws.Range("H16").Value = "(Romell) Niemelä, T.Wagner & M.Fisch."
$ws.Range("K16").Value = ""
$ws.Range("L16").Value = ""
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = ""
$ws.Range("Q16").Value = 425059
$ws.Range("R16").Value = 6712253
$ws.Range("AC16").Value = ""

# Row 17
$ws.Range("A17").Value = 131066768
$ws.Range("B17").Value = 91808
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 1202
$ws.Range("F17").Value = "Ullticka"
$ws.Range("G17").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H17").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K17").Value = ""
$ws.Range("L17").Value = ""
$ws.Range("M17").Value = ""
$ws.Range("N17").Value = ""
$ws.Range("Q17").Value = 425256
$ws.Range("R17").Value = 6712203
$ws.Range("AC17").Value = ""

# Row 18
$ws.Range("A18").Value = 131066769
$ws.Range("B18").Value = 57884
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 100109
$ws.Range("F18").Value = "Tretåig hackspett"
$ws.Range("G18").Value = "Picoides tridactylus"
$ws.Range("H18").Value = "(Linnaeus, 1758)"
$ws.Range("K18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = "äldre spår"
$ws.Range("N18").Value = ""
$ws.Range("Q18").Value = 425267
$ws.Range("R18").Value = 6712232
$ws.Range("AC18").Value = "Ringhack på gran"

# Row 30
$ws.Range("A30").Value = 131066776
$ws.Range("B30").Value = 80349
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 2081
$ws.Range("F30").Value = "Skrovellav"
$ws.Range("G30").Value = "Lobaria scrobiculata"
$ws.Range("H30").Value = "(Scop.) DC."
$ws.Range("K30").Value = ""
$ws.Range("L30").Value = ""
$ws.Range("M30").Value = ""
$ws.Range("N30").Value = ""
$ws.Range("Q30").Value = 425069
$ws.Range("R30").Value = 6712285
$ws.Range("AC30").Value = ""

# Row 31
$ws.Range("A31").Value = 131066772
$ws.Range("B31").Value = 57884
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 100109
$ws.Range("F31").Value = "Tretåig hackspett"
$ws.Range("G31").Value = "Picoides tridactylus"
$ws.Range("H31").Value = "(Linnaeus, 1758)"
$ws.Range("K31").Value = ""
$ws.Range("L31").Value = ""
$ws.Range("M31").Value = "äldre spår"
$ws.Range("N31").Value = ""
$ws.Range("Q31").Value = 425301
$ws.Range("R31").Value = 6712219
$ws.Range("AC31").Value = "Ringhack på gran"

# Row 32
$ws.Range("A32").Value = 131066767
$ws.Range("B32").Value = 91808
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 1202
$ws.Range("F32").Value = "Ullticka"
$ws.Range("G32").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H32").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K32").Value = ""
$ws.Range("L32").Value = ""
$ws.Range("M32").Value = ""
$ws.Range("N32").Value = ""
$ws.Range("Q32").Value = 425259
$ws.Range("R32").Value = 6712201
$ws.Range("AC32").Value = ""

# Row 33
$ws.Range("A33").Value = 131066790
$ws.Range("B33").Value = 83215
$ws.Range("D33").Value = "NT"
$ws.Range("E33").Value = 308
$ws.Range("F33").Value = "Brunpudrad nållav"
$ws.Range("G33").Value = "Chaenotheca gracillima"
$ws.Range("H33").Value = "(Vain.) Tibell"
$ws.Range("K33").Value = ""
$ws.Range("L33").Value = ""
$ws.Range("M33").Value = ""
$ws.Range("N33").Value = ""
$ws.Range("Q33").Value = 425164
$ws.Range("R33").Value = 6712278
$ws.Range("AC33").Value = ""
